$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ExtractedScans")

# Correct the Label values for the first three records from placeholder
# text values ("1001","1002","1003") to the real numeric scan IDs.
$ws.Range("B2").Value = 8097
$ws.Range("B3").Value = 8096
$ws.Range("B4").Value = 8095

$ws.Range("D10").Select()
